# Updated cryptos list on Thu Sep  7 17:30:46 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) and Volume(1h) (column E) figures
# for the crypto ranking table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.767.54'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.628.68'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("D4").Value = '''0.996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.79%  '
$ws.Range("D5").Value = '''214.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").Value = '''0.255'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").Value = '''0.0632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '1.853.86'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '1.626.98'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '''0.552'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '''62.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '25.755.61'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '''190.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''6.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.29%  '
$ws.Range("D24").Value = '''0.997'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '''142.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.00%  '
$ws.Range("D27").Value = '''0.123'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("D28").Value = '''6.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").Value = '''15.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '''0.0493'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("D32").Value = '''3.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").Value = '1.140.86'
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("D42").Value = '''2.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '''101.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").Value = '1.763.85'
$ws.Range("D47").Value = '''55.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  +7.69%  '
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").Value = '''7.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.54%  '
